$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend formatting for the new row 110 by copying row 109's format
$ws.Range("A109:B109").Copy() | Out-Null
$ws.Range("A110:B110").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data for rows 2-110: (Name, Quantity)
$data = @(
    @('ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы', 73440),
    @('Сб. Фитонефрол (Урологический сбор) 50г', 8991),
    @('Сб. Грудной №4 50г', 15877),
    @('Мать-и-мачеха листья 35г', 12441),
    @('Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г', 3361),
    @('Тысячелистник трава 50г', 8339),
    @('Зверобой трава 50г', 18744),
    @('Укроп пахучий плоды 50г', 35449),
    @('Полынь горькая трава 50г', 28126),
    @('Валериана корневища с корнями 50г', 14254),
    @('Чистотел трава 50г', 14350),
    @('Девясил корневища и корни 50г', 13118),
    @('Ромашка цветки вн 50г', 90690),
    @('Эвкалипт прутовидный листья 75г', 19020),
    @('Кукуруза столбики с рыльцами 40г', 17073),
    @('Бессмертник песчаный цветки 30г', 20942),
    @('Сб. Фитопектол №2 (Грудной сбор №2) 35г', 7871),
    @('Крушина кора 50г', 9313),
    @('Солодка корни 50г', 29338),
    @('Ноготки цветки 50г', 20492),
    @('Боярышник плоды 75г', 18355),
    @('Сенна листья 50г', 21295),
    @('Багульник болотный побеги 50г', 15033),
    @('Лен семена 100г', 68815),
    @('Пустырник трава 50г', 11618),
    @('Шалфей листья 50г', 33827),
    @('Эрва шерстистая трава 30г', 17068),
    @('Дуба кора 75г', 81300),
    @('Шиповник плоды низковитаминные 50г', 27014),
    @('Пижма цветки 75г', 18373),
    @('Липа цветки 35г', 20507),
    @('Мята перечная листья 50г', 24626),
    @('Можжевельник плоды 50г', 15298),
    @('Брусника листья 50г', 19188),
    @('Сб. Фитопектол №1 (Грудной сбор №1) 35г', 7312),
    @('Чага (березовый гриб) 50г', 26020),
    @('Толокнянка листья 50г', 9279),
    @('Чабрец трава 50г', 27850),
    @('Череда трава 50г', 16683),
    @('Ламинарии слоевища (морская капуста) 100г', 21738),
    @('Спорыш трава 50г', 17292),
    @('Аир корневища 75г', 11129),
    @('Алтей корни 75г', 7237),
    @('Крапива листья 50г', 26010),
    @('Рябина плоды 50г', 2800),
    @('Береза почки 50г', 26147),
    @('Подорожник большой листья 50г', 17075),
    @('Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г', 1850),
    @('Фп "ФармаЦветик® Фиточай для кормящих мам" 20х1,5 г', 1928),
    @('Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г', 3220),
    @('Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г', 4480),
    @('Фп Фиточай "Лактафитол" (БАД) 20х1,5 г', 14338),
    @('Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г', 3770),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г', 2810),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г', 6320),
    @('Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г', 2818),
    @('Фп Пижма цветки 20х1,5г', 2334),
    @('Фп Шалфей листья 20х1,5г', 66852),
    @('Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г', 11736),
    @('Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г', 3816),
    @('Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г', 2888),
    @('Фп Сб. Желудочный №3 20x2,0г', 10818),
    @('Фп Сб. Бруснивер 20x2,0г', 107395),
    @('Фп Череда трава 20х1,5г', 26765),
    @('Фп Сб. Элекасол 20x2,0г', 24174),
    @('Фп Крапива листья 20x1,5г', 33017),
    @('Фп Пустырник трава 20x1,5г', 19314),
    @('Фп Боярышник плоды 20х3,0г', 6788),
    @('Фп Чабрец трава 20x1,5 г', 47875),
    @('Фп Сб. Грудной №4 20x2,0г', 697284),
    @('Фп Ольха соплодия 20х1,5г', 1562),
    @('Фп Мелисса лекарственная трава 20x1,5г', 24154),
    @('Фп Шиповник плоды 20х2,0г', 29484),
    @('Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г', 3834),
    @('Фп Сб. Арфазетин-Э 20x2,0г', 30807),
    @('Фп Липа цветки 20x1,5г', 50814),
    @('Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г', 61274),
    @('Фп Хвощ полевой трава 20х1,5г', 22270),
    @('Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г', 60399),
    @('Фп Зверобой трава 20x1,5г', 34128),
    @('Фп Брусника листья 20х1,5г', 63841),
    @('Фп Пастушья сумка трава 20х1,5г', 4102),
    @('Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г', 151953),
    @('Фп Подорожник листья 20x1,5г', 23631),
    @('Фп Толокнянка листья 20x1,5г', 28322),
    @('Фп Чистотел трава 20х1,5г', 30663),
    @('Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г', 42460),
    @('Фп Ромашка цветки 20x1,5г', 1130381),
    @('Фп Душица трава 20x1,5г', 22819),
    @('Фп Мята перечная листья 20x1,5г', 53667),
    @('Фп Сенна листья 20x1,5г', 78515),
    @('Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г', 63018),
    @('Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г', 5040),
    @('Фп Береза листья 20x1,5г', 3741),
    @('Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г', 4212),
    @('Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г', 1836),
    @('Фп Золототысячник трава 20х1,5г', 7516),
    @('Фп Фиалка трехцветная трава 20x1,5г', 9054),
    @('Фп Аир корневища 20x1,5г', 11704),
    @('Фп Кровохлебка корневища и корни 20x1,5г', 3148),
    @('Фп Ноготки цветки 20x1,5г', 24386),
    @('Фп Крушина кора 20x1,5г', 7452),
    @('Фп Почечный чай листья 20x1,5г', 52040),
    @('Фп Валериана корневища с корнями 20x1,5г', 14622),
    @('Фп Тысячелистник трава 20x1,5г', 19386),
    @('Фп Лапчатка корневища 20x2,5г', 4050),
    @('Фп Дуб кора 20х1,5г', 7245),
    @('Фп Девясил корневища и корни 20х1,5г', 21626),
    @('Фп Бадан корневища 20x1,5г', 3733)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Update sheet view selection state
$ws.Activate()
$ws.Range("A88").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
